$d = $word.ActiveDocument

# Locate the target paragraph ("Entering this updates the fields below")
# and split its single run into three runs:
#   "Entering this " | "updates," | " the fields below"
$rng = $d.Content
$found = $rng.Find.Execute("Entering this updates the fields below", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
           '<w:body>' +
           '<w:p w14:paraId="63B5EDB0" w14:textId="51183587" w:rsidR="00214EF3" w:rsidRDefault="00214EF3" w:rsidP="00FE3CB4">' +
           '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
           '<w:r><w:t xml:space="preserve">Entering this </w:t></w:r>' +
           '<w:r><w:t>updates,</w:t></w:r>' +
           '<w:r><w:t xml:space="preserve"> the fields below</w:t></w:r>' +
           '</w:p>' +
           '</w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}
